# Applies the cryptos.xlsx update described in the commit diff.
# D-column price cells are numeric-looking text (e.g. "1.001", "0.06690");
# a leading apostrophe forces Excel to store them as text (quote-prefixed),
# preserving exact formatting (trailing zeros, thousand-dot grouping, etc.)
# exactly like the source inline strings.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'30.412.76"
$ws.Range("E2").Value = "  -0.85%  "

# Row 3
$ws.Range("D3").Value = "'1.890.08"
$ws.Range("E3").Value = "  +0.04%  "

# Row 4
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.03%  "

# Row 5
$ws.Range("D5").Value = "'237.92"
$ws.Range("E5").Value = "  +0.48%  "

# Row 6
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  +0.10%  "

# Row 7
$ws.Range("D7").Value = "'0.4893"
$ws.Range("E7").Value = "  +0.36%  "

# Row 8
$ws.Range("E8").Value = "  +0.28%  "

# Row 9
$ws.Range("D9").Value = "'0.06690"
$ws.Range("E9").Value = "  +0.22%  "

# Row 10
$ws.Range("D10").Value = "'1.887.17"
$ws.Range("E10").Value = "  -0.07%  "

# Row 11
$ws.Range("D11").Value = "'16.95"
$ws.Range("E11").Value = "  +1.47%  "

# Row 12
$ws.Range("D12").Value = "'0.07347"
$ws.Range("E12").Value = "  +1.54%  "

# Row 13
$ws.Range("D13").Value = "'5.124"
$ws.Range("E13").Value = "  +2.44%  "

# Row 14
$ws.Range("D14").Value = "'87.59"
$ws.Range("E14").Value = "  -1.85%  "

# Row 15
$ws.Range("D15").Value = "'0.6625"
$ws.Range("E15").Value = "  -0.22%  "

# Row 16
$ws.Range("D16").Value = "'30.389.85"
$ws.Range("E16").Value = "  -0.74%  "

# Row 17
$ws.Range("D17").Value = "'13.43"
$ws.Range("E17").Value = "  +3.19%  "

# Row 18
$ws.Range("D18").Value = "'0.000007821"
$ws.Range("E18").Value = "  -1.21%  "

# Row 20
$ws.Range("D20").Value = "'2.165.36"
$ws.Range("E20").Value = "  +1.39%  "

# Row 21
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "'5.304"
$ws.Range("E21").Value = "  +11.76%  "

# Row 22
$ws.Range("B22").Value = "BinanceUSD"
$ws.Range("C22").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D22").Value = "'1.000"
$ws.Range("E22").Value = "  -0.07%  "

# Row 23
$ws.Range("D23").Value = "'190.14"
$ws.Range("E23").Value = "  -1.42%  "

# Row 24
$ws.Range("D24").Value = "'6.107"
$ws.Range("E24").Value = "  +0.58%  "

# Row 25
$ws.Range("D25").Value = "'9.455"
$ws.Range("E25").Value = "  +1.63%  "

# Row 26
$ws.Range("D26").Value = "'163.27"
$ws.Range("E26").Value = "  +2.21%  "

# Row 27
$ws.Range("D27").Value = "'18.23"
$ws.Range("E27").Value = "  -0.36%  "

# Row 28
$ws.Range("D28").Value = "'1.927"
$ws.Range("E28").Value = "  +5.24%  "

# Row 29
$ws.Range("E29").Value = "  +4.69%  "

# Row 30
$ws.Range("D30").Value = "'4.355"
$ws.Range("E30").Value = "  +2.27%  "

# Row 31
$ws.Range("D31").Value = "'0.09144"
$ws.Range("E31").Value = "  +1.40%  "

# Row 32
$ws.Range("D32").Value = "'4.032"
$ws.Range("E32").Value = "  +2.40%  "

# Row 33
$ws.Range("D33").Value = "'0.05196"
$ws.Range("E33").Value = "  -0.06%  "

# Row 34
$ws.Range("D34").Value = "'0.7396"
$ws.Range("E34").Value = "  +0.89%  "

# Row 35
$ws.Range("D35").Value = "'1.097"
$ws.Range("E35").Value = "  +1.06%  "

# Row 36
$ws.Range("E36").Value = "  +1.37%  "

# Row 37
$ws.Range("D37").Value = "'0.01812"
$ws.Range("E37").Value = "  -0.57%  "

# Row 38
$ws.Range("D38").Value = "'2.664"
$ws.Range("E38").Value = "  -0.18%  "

# Row 39
$ws.Range("D39").Value = "'0.9210"
$ws.Range("E39").Value = "  -0.41%  "

# Row 40
$ws.Range("D40").Value = "'2.035"
$ws.Range("E40").Value = "  -0.66%  "

# Row 42
$ws.Range("D42").Value = "'5.930"
$ws.Range("E42").Value = "  +3.51%  "

# Row 43
$ws.Range("D43").Value = "'106.16"
$ws.Range("E43").Value = "  +1.52%  "

# Row 44
$ws.Range("D44").Value = "'0.9917"
$ws.Range("E44").Value = "  -0.76%  "

# Row 45
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Value = "'68.32"
$ws.Range("E45").Value = "  +19.38%  "

# Row 46
$ws.Range("B46").Value = "Algorand"
$ws.Range("C46").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D46").Value = "'0.1368"
$ws.Range("E46").Value = "  +2.35%  "

# Row 47
$ws.Range("D47").Value = "'7.554"
$ws.Range("E47").Value = "  +2.94%  "

# Row 48
$ws.Range("B48").Value = "Elrond"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D48").Value = "'34.92"
$ws.Range("E48").Value = "  +5.08%  "

# Row 49
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'8.936"
$ws.Range("E49").Value = "  +2.38%  "

# Row 50
$ws.Range("D50").Value = "'0.05820"
$ws.Range("E50").Value = "  -0.20%  "

# Row 51
$ws.Range("D51").Value = "'0.3936"
$ws.Range("E51").Value = "  -4.89%  "
